$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new clock-in/clock-out entry for 2026-01-23 ---
# Force column A to be treated as text first so the date-like string isn't
# auto-converted into a date serial number (it should stay plain text like
# the other date cells in the sheet).
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2026-01-23"
$ws.Range("B6").Value = "15:48:02"
$ws.Range("C6").Value = "15:48:11"
$ws.Range("D6").Value = "0 Hours"

# Re-apply the same formatting as the row above (row 5) to the new row 6,
# so it keeps the same cell style used throughout the table.
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 7: "Total Duration" summary moves down from row 5 to row 7 ---
$ws.Range("C5:D5").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C7").Value = "Total Duration:"
$ws.Range("D7").Value = "-10.5 Hours"

# The old C5/D5 "Total Duration:" values are removed; those cells become
# blank (numeric/empty) now that the summary lives on row 7.
$ws.Range("C5:D5").ClearContents()

$excel.CutCopyMode = $false
